# Remove stray apostrophes from facility names in column B
# ("remove dupes" - de-duplicating facility-name variants that differ
# only by an apostrophe/possessive form).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$col = $ws.Columns.Item(2)

$col.Replace("Olin E. Teague Veterans' Center", "Olin E. Teague Veterans Center", -4142, 1, $false, $false, $true)
$col.Replace("Audie L. Murphy Memorial Veterans' Hospital", "Audie L. Murphy Memorial Veterans Hospital", -4142, 1, $false, $false, $true)
$col.Replace("George H. O'Brien, Jr. Department of Veterans Affairs Medical Center", "George H. OBrien, Jr. Department of Veterans Affairs Medical Center", -4142, 1, $false, $false, $true)
